# fix(publipostage): Refactor synthetic array /3
# Swap the "noir" (black) status-color entry for a "bleu" (blue) one:
#   square emoji ⬛/🟥/🟧/🟩 -> book emoji 📘/📕/📙/📗
#   label "noir" -> "bleu" (rouge/orange/vert stay as-is)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emojiMap = @{
    "⬛" = "📘"
    "🟥" = "📕"
    "🟧" = "📙"
    "🟩" = "📗"
}

$labelMap = @{
    "noir" = "bleu"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$changedA = 0
$changedB = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $txtA = $cellA.Text
    if ($emojiMap.ContainsKey($txtA)) {
        $cellA.Value = $emojiMap[$txtA]
        $changedA = $changedA + 1
    }

    $cellB = $ws.Cells.Item($r, 2)
    $txtB = $cellB.Text
    if ($labelMap.ContainsKey($txtB)) {
        $cellB.Value = $labelMap[$txtB]
        $changedB = $changedB + 1
    }
}

Write-Output "Changed A: $changedA, Changed B: $changedB"
